$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F ("Is Active") currently stores the text strings "True"/"False".
# Push up real boolean values instead, so Excel writes them as native
# boolean cells (t="b") rather than shared-string text.
$ws.Range("F2").Value = $true
$ws.Range("F3").Value = $true
$ws.Range("F4").Value = $false
$ws.Range("F5").Value = $true
